$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 41 (pushing existing rows 41-67 down to 43-69)
$ws.Rows.Item(41).Insert()
$ws.Rows.Item(41).Insert()

# New row 41 data
$ws.Cells.Item(41, 1).Value = 6
$ws.Cells.Item(41, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(41, 3).Value = "Metropolitana"
$ws.Cells.Item(41, 4).Value = 44651
$ws.Cells.Item(41, 5).Value = 13
$ws.Cells.Item(41, 6).Value = "Fruta"
$ws.Cells.Item(41, 7).Value = 100104
$ws.Cells.Item(41, 8).Value = "Frutos de pepita"
$ws.Cells.Item(41, 9).Value = 100104003
$ws.Cells.Item(41, 10).Value = "Membrillo"
$ws.Cells.Item(41, 11).Value = "Champion"
$ws.Cells.Item(41, 12).Value = "Especial"
$ws.Cells.Item(41, 13).Value = 8
$ws.Cells.Item(41, 14).Value = 270000
$ws.Cells.Item(41, 15).Value = 270000
$ws.Cells.Item(41, 16).Value = 270000
$ws.Cells.Item(41, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(41, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(41, 19).Value = 600
$ws.Cells.Item(41, 20).Value = 450

# New row 42 data
$ws.Cells.Item(42, 1).Value = 6
$ws.Cells.Item(42, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(42, 3).Value = "Metropolitana"
$ws.Cells.Item(42, 4).Value = 44651
$ws.Cells.Item(42, 5).Value = 13
$ws.Cells.Item(42, 6).Value = "Fruta"
$ws.Cells.Item(42, 7).Value = 100104
$ws.Cells.Item(42, 8).Value = "Frutos de pepita"
$ws.Cells.Item(42, 9).Value = 100104003
$ws.Cells.Item(42, 10).Value = "Membrillo"
$ws.Cells.Item(42, 11).Value = "Champion"
$ws.Cells.Item(42, 12).Value = "Primera"
$ws.Cells.Item(42, 13).Value = 27
$ws.Cells.Item(42, 14).Value = 230000
$ws.Cells.Item(42, 15).Value = 250000
$ws.Cells.Item(42, 16).Value = 238889
$ws.Cells.Item(42, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(42, 18).Value = "Región Metropolitana"
$ws.Cells.Item(42, 19).Value = 531
$ws.Cells.Item(42, 20).Value = 450
